$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps the same text "HK_G_acc_SD" but the diff shows a new shared
# string entry was introduced (a duplicate) and A1 now points at it.
# Re-assigning the text value reproduces the same observable state.
$ws.Range("A1").Value = "HK_G_acc_SD"

$values = @(
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.666666666666671,
    86.666666666666671,
    86.666666666666671,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    87.567567567567579,
    87.387387387387378,
    87.567567567567579,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    87.567567567567579,
    87.747747747747752,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    86.486486486486484,
    87.747747747747752,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    87.387387387387378,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297,
    86.306306306306297
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
